$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SECTION (column B) values from "B" to "A" for the affected roll numbers
$ws.Range("B3").Value = "A"
$ws.Range("B4").Value = "A"
$ws.Range("B6").Value = "A"
$ws.Range("B44").Value = "A"
$ws.Range("B45").Value = "A"
$ws.Range("B47").Value = "A"
$ws.Range("B49").Value = "A"

# Update the active cell selection to C2
$ws.Range("C2").Select()
